$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The worksheet is protected; unprotect it so the cell values can be updated
$ws.Unprotect()

# Update the confidential disclaimer text in cell A59 (shared string), changing the date from 2021-05-20 to 2021-05-21
$oldText = $ws.Range("A59").Value2
$newText = $oldText -replace [regex]::Escape("2021-05-20"), "2021-05-21"
$ws.Range("A59").Value = $newText

# Update the Weight (column D) and Percent Change (column E) values for rows 2-56
$ws.Range("D2").Value = 0.01441195565244332
$ws.Range("E2").Value = 0.0060690943043884
$ws.Range("D3").Value = 0.05058164329517143
$ws.Range("E3").Value = -0.01373288008670803
$ws.Range("D4").Value = 0.01449350478422155
$ws.Range("E4").Value = -0.006340131960712814
$ws.Range("D5").Value = 0.00977094410091473
$ws.Range("E5").Value = 0.006885998469777999
$ws.Range("D6").Value = 0.01559282944377741
$ws.Range("E6").Value = -0.0008829722203355095
$ws.Range("D7").Value = 0.02020350147874679
$ws.Range("E7").Value = 0
$ws.Range("D8").Value = 0.00459749583835281
$ws.Range("E8").Value = 0.009058572444865964
$ws.Range("D9").Value = 0.006598190714206076
$ws.Range("E9").Value = 0.001359619306593896
$ws.Range("D10").Value = 0.0143884067052988
$ws.Range("E10").Value = 0.0002532928064842288
$ws.Range("D11").Value = 0.008265805045328783
$ws.Range("E11").Value = -0.002487187217365183
$ws.Range("D12").Value = 0.01552361547478916
$ws.Range("E12").Value = 0.009390801830002404
$ws.Range("D13").Value = 0.003001400531891223
$ws.Range("E13").Value = -0.01307664366146022
$ws.Range("D14").Value = 0.005814658681834201
$ws.Range("E14").Value = -0.001767825574543136
$ws.Range("D15").Value = 0.01452829866512158
$ws.Range("E15").Value = 0.0113784741652676
$ws.Range("D16").Value = 0.01068112959769083
$ws.Range("E16").Value = 0.01902887139107623
$ws.Range("D17").Value = 0.02133142128840671
$ws.Range("E17").Value = -0.005256932579839835
$ws.Range("D18").Value = 0.008765441437125111
$ws.Range("E18").Value = 0.0007462686567163423
$ws.Range("D19").Value = 0.0170519296847078
$ws.Range("E19").Value = -0.0006430116326648916
$ws.Range("D20").Value = 0.01243047995699768
$ws.Range("E20").Value = 0.001754122187139773
$ws.Range("D21").Value = 0.006832185014403956
$ws.Range("E21").Value = 0.01248313090418351
$ws.Range("D22").Value = 0.0155370720160146
$ws.Range("E22").Value = 0.008308072302683378
$ws.Range("D23").Value = 0.0195153489121904
$ws.Range("E23").Value = -0.0002154800896396747
$ws.Range("D24").Value = 0.009918592261582721
$ws.Range("E24").Value = 0.006375227686703067
$ws.Range("D25").Value = 0.01994907317150288
$ws.Range("E25").Value = 0.003572588502760654
$ws.Range("D26").Value = 0.01399461597804885
$ws.Range("E26").Value = -0.009675609984107547
$ws.Range("D27").Value = 0.0214741478437191
$ws.Range("E27").Value = -0.02013658412050034
$ws.Range("D28").Value = 0.05551882335126358
$ws.Range("E28").Value = -0.01476710391956648
$ws.Range("D29").Value = 0.02055455522784558
$ws.Range("E29").Value = -0.003890160183066538
$ws.Range("D30").Value = 0.02982617443097383
$ws.Range("E30").Value = -0.005054724705489289
$ws.Range("D31").Value = 0.01528712922251249
$ws.Range("E31").Value = -0.008896260554885238
$ws.Range("D32").Value = 0.01318734810212656
$ws.Range("E32").Value = -0.004398168925590085
$ws.Range("D33").Value = 0.01819548649366209
$ws.Range("E33").Value = -0.005053617650684017
$ws.Range("D34").Value = 0.04311606636114237
$ws.Range("E34").Value = -0.005557120873880939
$ws.Range("D35").Value = 0.01096708109873138
$ws.Range("E35").Value = 0.005453306066802943
$ws.Range("D36").Value = 0.01021413777812658
$ws.Range("E36").Value = -0.008453590641277375
$ws.Range("D37").Value = 0.01060870724040115
$ws.Range("E37").Value = -0.004404316229905225
$ws.Range("D38").Value = 0.007530056194067606
$ws.Range("E38").Value = 0.001365103003226453
$ws.Range("D39").Value = 0.01212929639887557
$ws.Range("E39").Value = 0.01265822784810133
$ws.Range("D40").Value = 0.01745020592569957
$ws.Range("E40").Value = -0.007947019867549754
$ws.Range("D41").Value = 0.01700115616110257
$ws.Range("E41").Value = -0.001392466754856159
$ws.Range("D42").Value = 0.03271341240825573
$ws.Range("E42").Value = -0.003456451566829388
$ws.Range("D43").Value = 0.01139071544407685
$ws.Range("E43").Value = 0.001957776863317218
$ws.Range("D44").Value = 0.02186575811290011
$ws.Range("E44").Value = 0.001457339692633841
$ws.Range("D45").Value = 0.01275567970327908
$ws.Range("E45").Value = -0.008263736263736332
$ws.Range("D46").Value = 0.00866405013691956
$ws.Range("E46").Value = 0.001779647306261012
$ws.Range("D47").Value = 0.01317563592735627
$ws.Range("E47").Value = 0.009574875526618065
$ws.Range("D48").Value = 0.01032789539052311
$ws.Range("E48").Value = 0.01653094462540716
$ws.Range("D49").Value = 0.015790503542427
$ws.Range("E49").Value = 0.01880740460183694
$ws.Range("D50").Value = 0.008561599757080525
$ws.Range("E50").Value = -0.00185551707075704
$ws.Range("D51").Value = 0.01083183039965534
$ws.Range("E51").Value = -0.01070921210796638
$ws.Range("D52").Value = 0.0083875680537784
$ws.Range("E52").Value = -0.001266391600964201
$ws.Range("D53").Value = 0.009588564358148675
$ws.Range("E53").Value = 0.01421262080727681
$ws.Range("D54").Value = 0.1353192277582051
$ws.Range("D55").Value = 0.04378764744637485
$ws.Range("E55").Value = -0.0009674688595959768
$ws.Range("D56").Value = 0.9999999999999999
$ws.Range("E56").Value = -0.001716490233220069

# Restore the worksheet's protected state (original password is not recoverable from the
# stored legacy hash, so we re-protect without a password to preserve the sheet's protected status)
$ws.Protect()
